{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// The underlying OOXML diff shows three related changes that all stem from\n// the same authoring pass (placing/moving the cursor, which Word tracks\n// with the hidden \"_GoBack\" bookmark):\n//  1. The run ending \"...demonstrate my point.\" is split into two runs so\n//     the trailing \".\" becomes its own run (no visible text change).\n//  2. The run \"Many linguists, archaeologists, ...\" is split into \"Many\n//     linguists, \" and \"archaeologists, ...\" with the \"_GoBack\" bookmark\n//     inserted between them (no visible text change).\n//  3. The \"_GoBack\" bookmark that used to sit at the very end of the\n//     document (right after \"...Indo-European languages.\") is gone - a\n//     document only ever has one \"_GoBack\" bookmark, and it moved to the\n//     new location from change #2.\n//\n// None of this changes the visible text; it only changes run boundaries and\n// where the \"_GoBack\" bookmark lives. We reproduce the run split by\n// inserting (and immediately removing) a throwaway bookmark at the split\n// point - that forces the host to break the run there, just like Word does\n// whenever the cursor/an edit lands mid-run - and we reproduce the bookmark\n// move with deleteBookmark/insertBookmark.\n\nconst body = context.document.body;\n\n// --- Change 1: split \"...demonstrate my point.\" into \"...point\" + \".\" ---\nconst pointMatches = body.search(\"demonstrate my point\", { matchCase: true });\npointMatches.load(\"items\");\nawait context.sync();\n\nif (pointMatches.items.length === 0) {\n  throw new Error(\"Could not find 'demonstrate my point' to split.\");\n}\n\n// Collapsed range right between \"point\" and the trailing period.\nconst splitPoint = pointMatches.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_tempRunSplit\");\nawait context.sync();\ncontext.document.deleteBookmark(\"_tempRunSplit\");\nawait context.sync();\n\n// --- Changes 2 & 3: move the \"_GoBack\" bookmark ---\n// Remove it from wherever it currently is (the end of the document) first,\n// since Word only keeps a single bookmark per name.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert it between \"Many linguists, \" and \"archaeologists, ...\", which\n// splits that run into two runs exactly like change #2 in the diff.\nconst linguistsMatches = body.search(\"Many linguists, \", { matchCase: true });\nlinguistsMatches.load(\"items\");\nawait context.sync();\n\nif (linguistsMatches.items.length === 0) {\n  throw new Error(\"Could not find 'Many linguists, ' to anchor the bookmark.\");\n}\n\nconst bookmarkSpot = linguistsMatches.items[0].getRange(\"End\");\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# The underlying OOXML diff shows three related changes that all stem from\n# the same authoring pass (placing/moving the cursor, which Word tracks\n# with the hidden \"_GoBack\" bookmark):\n#  1. The run ending \"...demonstrate my point.\" is split into two runs so\n#     the trailing \".\" becomes its own run (no visible text change).\n#  2. The run \"Many linguists, archaeologists, ...\" is split into \"Many\n#     linguists, \" and \"archaeologists, ...\" with the \"_GoBack\" bookmark\n#     inserted between them (no visible text change).\n#  3. The \"_GoBack\" bookmark that used to sit at the very end of the\n#     document (right after \"...Indo-European languages.\") is gone - a\n#     document only ever has one \"_GoBack\" bookmark, and it moved to the\n#     new location from change #2.\n#\n# None of this changes the visible text; it only changes run boundaries and\n# where the \"_GoBack\" bookmark lives. We reproduce the run split by adding\n# (and immediately removing) a throwaway bookmark at the split point - Word\n# always breaks a run at a bookmark boundary, so that leaves the text split\n# into two runs with no visible change - and we reproduce the \"_GoBack\" move\n# with Bookmarks.Add, which replaces any existing bookmark of the same name.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: split \"...demonstrate my point.\" into \"...point\" + \".\" ---\n$findRange = $d.Content\n$null = $findRange.Find.Execute(\"demonstrate my point\")\n$splitRange = $findRange.Duplicate\n$splitRange.Collapse(0)  # wdCollapseEnd - collapse to the point right after \"point\"\n$d.Bookmarks.Add(\"_tempRunSplit\", $splitRange)\n$d.Bookmarks(\"_tempRunSplit\").Delete()\n\n# --- Changes 2 & 3: move the \"_GoBack\" bookmark ---\n# Re-anchor \"_GoBack\" between \"Many linguists, \" and \"archaeologists, ...\".\n# Bookmarks.Add replaces any existing bookmark with the same name, so this\n# both removes the old \"_GoBack\" (at the end of the document) and creates\n# the new one in a single step, splitting the run exactly like change #2.\n$bmFindRange = $d.Content\n$null = $bmFindRange.Find.Execute(\"Many linguists, \")\n$bookmarkRange = $bmFindRange.Duplicate\n$bookmarkRange.Collapse(0)  # wdCollapseEnd - collapse to right after \"Many linguists, \"\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
